$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 4 for columns A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr4 = "${col}4"
    $val2 = $ws.Range($addr2).Value()
    $val4 = $ws.Range($addr4).Value()
    $ws.Range($addr2).Value = $val4
    $ws.Range($addr4).Value = $val2
}
